$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.407.71"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "2.929.92"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'593.67"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "'143.54"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.499"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").Value = "'6.94"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").Value = "'33.18"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").Value = "3.417.95"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "61.387.95"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "2.930.98"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "'6.62"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "'433.09"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "'13.53"
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "'7.05"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "'81.43"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'10.83"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").Value = "'11.71"
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -3.65%  "
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("D31").Value = "'26.71"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  +2.68%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "'5.60"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").Value = "'1.98"
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").Value = "'41.91"
$ws.Range("E41").Value = "  +5.06%  "
$ws.Range("E42").Value = "  -2.32%  "
$ws.Range("D43").Value = "'0.0343"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "2.694.36"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "'133.12"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").Value = "'363.34"
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "'23.51"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").Value = "'1.99"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("E51").Value = "  +0.53%  "
